$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 11.00503465770099
$ws.Range("C2").Value = 9.860640405116646
$ws.Range("D2").Value = 5.883248273967437
$ws.Range("F2").Value = 27.97989136244867
$ws.Range("G2").Value = 3.660796443117949
$ws.Range("I2").Value = 25.70650928085748
$ws.Range("K2").Value = 8.140289405652602
$ws.Range("L2").Value = 11.17229863751233
$ws.Range("M2").Value = 14.19007656298228
$ws.Range("O2").Value = 25.14335259637067
$ws.Range("B3").Value = 10.75536320994728
$ws.Range("C3").Value = 9.846550413190339
$ws.Range("D3").Value = 5.833088253940462
$ws.Range("F3").Value = 28.03194421124823
$ws.Range("G3").Value = 3.6625355580594
$ws.Range("I3").Value = 25.80110696411576
$ws.Range("K3").Value = 7.948690004046215
$ws.Range("L3").Value = 11.1815379742186
$ws.Range("M3").Value = 14.15337675779295
$ws.Range("O3").Value = 25.22329460217201
$ws.Range("B4").Value = 10.6008167639836
$ws.Range("C4").Value = 9.837998286955269
$ws.Range("D4").Value = 5.801571884076187
$ws.Range("F4").Value = 28.0705530818498
$ws.Range("G4").Value = 3.663660554525146
$ws.Range("I4").Value = 25.86365266166924
$ws.Range("K4").Value = 7.82971824881734
$ws.Range("L4").Value = 11.18892913295986
$ws.Range("M4").Value = 14.13287277290985
$ws.Range("O4").Value = 25.27723960712031
$ws.Range("B5").Value = 10.53761492041673
$ws.Range("C5").Value = 9.83453906704924
$ws.Range("D5").Value = 5.788553091610778
$ws.Range("F5").Value = 28.0879553250115
$ws.Range("G5").Value = 3.66413342036395
$ws.Range("I5").Value = 25.89026221378755
$ws.Range("K5").Value = 7.780970590208806
$ws.Range("L5").Value = 11.19237359789029
$ws.Range("M5").Value = 14.12503361667828
$ws.Range("O5").Value = 25.30044318271412
$ws.Range("B6").Value = 10.52710955437921
$ws.Range("C6").Value = 9.833966253087299
$ws.Range("D6").Value = 5.786380868912155
$ws.Range("F6").Value = 28.09094565117763
$ws.Range("G6").Value = 3.664212811644187
$ws.Range("I6").Value = 25.89474844392193
$ws.Range("K6").Value = 7.772862088976937
$ws.Range("L6").Value = 11.1929716825391
$ws.Range("M6").Value = 14.1237632852976
$ws.Range("O6").Value = 25.304369788568
$ws.Range("B7").Value = 10.59996518209634
$ws.Range("C7").Value = 9.837951528892837
$ws.Range("D7").Value = 5.801397012629305
$ws.Range("F7").Value = 28.07078102187032
$ws.Range("G7").Value = 3.663666873312579
$ws.Range("I7").Value = 25.86400698633827
$ws.Range("K7").Value = 7.82906180404007
$ws.Range("L7").Value = 11.18897383447514
$ws.Range("M7").Value = 14.13276495291381
$ws.Range("O7").Value = 25.27754759888879
$ws.Range("B8").Value = 10.91926013419522
$ws.Range("C8").Value = 9.855762129294206
$ws.Range("D8").Value = 5.866105023527327
$ws.Range("F8").Value = 27.99645772941896
$ws.Range("G8").Value = 3.661384250289538
$ws.Range("I8").Value = 25.73820015378546
$ws.Range("K8").Value = 8.074541899372617
$ws.Range("L8").Value = 11.17512805891214
$ws.Range("M8").Value = 14.17700487197114
$ws.Range("O8").Value = 25.16990689755172
$ws.Range("B9").Value = 11.53170222326224
$ws.Range("C9").Value = 9.891437079451872
$ws.Range("D9").Value = 5.987090555901492
$ws.Range("F9").Value = 27.90357224007345
$ws.Range("G9").Value = 3.657359663086222
$ws.Range("I9").Value = 25.5269184551805
$ws.Range("K9").Value = 8.542535121834231
$ws.Range("L9").Value = 11.16158617989361
$ws.Range("M9").Value = 14.27958681603031
$ws.Range("O9").Value = 24.99745804377318
$ws.Range("B10").Value = 11.96852682450382
$ws.Range("C10").Value = 9.918053342599071
$ws.Range("D10").Value = 6.072104077648014
$ws.Range("F10").Value = 27.86767906342812
$ws.Range("G10").Value = 3.654675307896885
$ws.Range("I10").Value = 25.39331081014716
$ws.Range("K10").Value = 8.874640639165385
$ws.Range("L10").Value = 11.15989796724426
$ws.Range("M10").Value = 14.36420667393092
$ws.Range("O10").Value = 24.89440495533887
$ws.Range("B11").Value = 12.16347822580899
$ws.Range("C11").Value = 9.930239078816019
$ws.Range("D11").Value = 6.109881366526182
$ws.Range("F11").Value = 27.85838979883902
$ws.Range("G11").Value = 3.653512699222399
$ws.Range("I11").Value = 25.3372296905568
$ws.Range("K11").Value = 9.022499389160169
$ws.Range("L11").Value = 11.16091394207632
$ws.Range("M11").Value = 14.4046194864403
$ws.Range("O11").Value = 24.85267620596873
$ws.Range("B12").Value = 12.23669048113209
$ws.Range("C12").Value = 9.934863723247192
$ws.Range("D12").Value = 6.124053168886225
$ws.Range("F12").Value = 27.8558845312153
$ws.Range("G12").Value = 3.653080819111362
$ws.Range("I12").Value = 25.31666933413077
$ws.Range("K12").Value = 9.077975851162147
$ws.Range("L12").Value = 11.1615541468363
$ws.Range("M12").Value = 14.42019030506381
$ws.Range("O12").Value = 24.83761668570327
$ws.Range("B13").Value = 12.22095118222156
$ws.Range("C13").Value = 9.933867287423798
$ws.Range("D13").Value = 6.121007036789621
$ws.Range("F13").Value = 27.85637906299462
$ws.Range("G13").Value = 3.653173460396522
$ws.Range("I13").Value = 25.32106728081778
$ws.Range("K13").Value = 9.066051662941103
$ws.Range("L13").Value = 11.16140492291256
$ws.Range("M13").Value = 14.41682509939343
$ws.Range("O13").Value = 24.84082698654918
$ws.Range("B14").Value = 12.16951408554181
$ws.Range("C14").Value = 9.930619353173915
$ws.Range("D14").Value = 6.11104999473775
$ws.Range("F14").Value = 27.85816340044632
$ws.Range("G14").Value = 3.653477000576048
$ws.Range("I14").Value = 25.33552461790474
$ws.Range("K14").Value = 9.027074063115109
$ws.Range("L14").Value = 11.16096149842072
$ws.Range("M14").Value = 14.40589519889284
$ws.Range("O14").Value = 24.85142236552363
$ws.Range("B15").Value = 12.13792572267377
$ws.Range("C15").Value = 9.928631193558338
$ws.Range("D15").Value = 6.104933474304623
$ws.Range("F15").Value = 27.85938819552178
$ws.Range("G15").Value = 3.653664017127807
$ws.Range("I15").Value = 25.34446825926942
$ws.Range("K15").Value = 9.003130696508556
$ws.Range("L15").Value = 11.16072312549129
$ws.Range("M15").Value = 14.39923487820156
$ws.Range("O15").Value = 24.85800904998752
$ws.Range("B16").Value = 11.95570443142657
$ws.Range("C16").Value = 9.917258455404038
$ws.Range("D16").Value = 6.069616802448202
$ws.Range("F16").Value = 27.86842780847893
$ws.Range("G16").Value = 3.654752461190897
$ws.Range("I16").Value = 25.39707044286359
$ws.Range("K16").Value = 8.864908458559242
$ws.Range("L16").Value = 11.15986738128918
$ws.Range("M16").Value = 14.36160348897396
$ws.Range("O16").Value = 24.89723580964861
$ws.Range("B17").Value = 11.84290207750149
$ws.Range("C17").Value = 9.910300864749368
$ws.Range("D17").Value = 6.047718266982384
$ws.Range("F17").Value = 27.87577645488735
$ws.Range("G17").Value = 3.655435145475921
$ws.Range("I17").Value = 25.4305439506291
$ws.Range("K17").Value = 8.779251362363469
$ws.Range("L17").Value = 11.15979861148226
$ws.Range("M17").Value = 14.339003321795
$ws.Range("O17").Value = 24.92262051991119
$ws.Range("B18").Value = 11.77767091356112
$ws.Range("C18").Value = 9.906306433609688
$ws.Range("D18").Value = 6.035038774182333
$ws.Range("F18").Value = 27.88066574601737
$ws.Range("G18").Value = 3.655833317793363
$ws.Range("I18").Value = 25.45023914369836
$ws.Range("K18").Value = 8.729683634715
$ws.Range("L18").Value = 11.1599270414018
$ws.Range("M18").Value = 14.32618544541974
$ws.Range("O18").Value = 24.93770584442036
$ws.Range("B19").Value = 11.75552684556627
$ws.Range("C19").Value = 9.904955288379577
$ws.Range("D19").Value = 6.030731431653118
$ws.Range("F19").Value = 27.88243495180489
$ws.Range("G19").Value = 3.655969079742209
$ws.Range("I19").Value = 25.45698350656648
$ws.Range("K19").Value = 8.712850965893995
$ws.Range("L19").Value = 11.15999940718532
$ws.Range("M19").Value = 14.32187689935254
$ws.Range("O19").Value = 24.94289667538514
$ws.Range("B20").Value = 11.85494685637479
$ws.Range("C20").Value = 9.911040754618975
$ws.Range("D20").Value = 6.05005813601173
$ws.Range("F20").Value = 27.8749256045506
$ws.Range("G20").Value = 3.655361902557758
$ws.Range("I20").Value = 25.42693487517814
$ws.Range("K20").Value = 8.788401146213003
$ws.Range("L20").Value = 11.15978855363057
$ws.Range("M20").Value = 14.34139046159983
$ws.Range("O20").Value = 24.91986809655315
$ws.Range("B21").Value = 12.18463954339881
$ws.Range("C21").Value = 9.93157308239241
$ws.Range("D21").Value = 6.113978284087035
$ws.Range("F21").Value = 27.85761182309563
$ws.Range("G21").Value = 3.653387616503685
$ws.Range("I21").Value = 25.33125978507288
$ws.Range("K21").Value = 9.038537082933969
$ws.Range("L21").Value = 11.16108481835122
$ws.Range("M21").Value = 14.40909838929081
$ws.Range("O21").Value = 24.84829008911382
$ws.Range("B22").Value = 12.39651876703996
$ws.Range("C22").Value = 9.945050944928372
$ws.Range("D22").Value = 6.15497270293757
$ws.Range("F22").Value = 27.85219681946088
$ws.Range("G22").Value = 3.652146101354594
$ws.Range("I22").Value = 25.27267323631574
$ws.Range("K22").Value = 9.198995687116478
$ws.Range("L22").Value = 11.16342054552204
$ws.Range("M22").Value = 14.45490363602067
$ws.Range("O22").Value = 24.80583680758674
$ws.Range("B23").Value = 12.28378566848667
$ws.Range("C23").Value = 9.937852517322135
$ws.Range("D23").Value = 6.133166243340678
$ws.Range("F23").Value = 27.85454709888645
$ws.Range("G23").Value = 3.652804269827881
$ws.Range("I23").Value = 25.30358095176925
$ws.Range("K23").Value = 9.113648260420216
$ws.Range("L23").Value = 11.16203810095448
$ws.Range("M23").Value = 14.43031721589288
$ws.Range("O23").Value = 24.82809849070447
$ws.Range("B24").Value = 11.84950258875958
$ws.Range("C24").Value = 9.910706232952936
$ws.Range("D24").Value = 6.049000560352378
$ws.Range("F24").Value = 27.87530820427246
$ws.Range("G24").Value = 3.655394997960457
$ws.Range("I24").Value = 25.42856513394054
$ws.Range("K24").Value = 8.784265528923195
$ws.Range("L24").Value = 11.15979257751166
$ws.Range("M24").Value = 14.34031068913115
$ws.Range("O24").Value = 24.92111093684372
$ws.Range("B25").Value = 11.36799880040075
$ws.Range("C25").Value = 9.881711096302613
$ws.Range("D25").Value = 5.95502169444312
$ws.Range("F25").Value = 27.92302597995842
$ws.Range("G25").Value = 3.658400362899957
$ws.Range("I25").Value = 25.58028072191957
$ws.Range("K25").Value = 8.4177470043266
$ws.Range("L25").Value = 11.16379586087651
$ws.Range("M25").Value = 14.27958681603031
$ws.Range("O25").Value = 25.03996410723348
